$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F-column "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 189
$ws1.Range("F3").Value = 5357
$ws1.Range("F4").Value = 30
$ws1.Range("F6").Value = 24
$ws1.Range("F7").Value = 606
$ws1.Range("F8").Value = 575
$ws1.Range("F12").Value = 4291
$ws1.Range("F13").Value = 439
$ws1.Range("F14").Value = 193
$ws1.Range("F15").Value = 166
$ws1.Range("F17").Value = 3423
$ws1.Range("F18").Value = 169
$ws1.Range("F19").Value = 1092
$ws1.Range("F22").Value = 200
$ws1.Range("F25").Value = 141
$ws1.Range("F27").Value = 309
$ws1.Range("F28").Value = 31

# Sheet "全部类型" (sheet4): F-column "想去人数" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 189
$ws4.Range("F4").Value = 5357
$ws4.Range("F5").Value = 30
$ws4.Range("F7").Value = 24
$ws4.Range("F8").Value = 606
$ws4.Range("F9").Value = 575
$ws4.Range("F13").Value = 4291
$ws4.Range("F14").Value = 439
$ws4.Range("F15").Value = 193
$ws4.Range("F16").Value = 166
$ws4.Range("F18").Value = 3423
$ws4.Range("F19").Value = 169
$ws4.Range("F20").Value = 1092
$ws4.Range("F23").Value = 200
$ws4.Range("F26").Value = 141
$ws4.Range("F28").Value = 309
$ws4.Range("F29").Value = 31
